$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.790271878242493
$ws.Range("B1").Value = 4.180408000946045
$ws.Range("C1").Value = 1.949701070785522
$ws.Range("D1").Value = 0.8850273489952087
$ws.Range("E1").Value = 0.4776319265365601
